$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append, replicating the "Quarterly Payout" pattern already
# present in the sheet, but with maturity dates moved to the 30th of the
# respective month.
$newRows = @(
    @("Quarterly Payout", "100000",  "5 Years: 0 Months : 0 Days",   "1826", "6.5",  "30 Apr 2029", "100000",  "₹32,500",  "1,625"),
    @("Quarterly Payout", "4000",    "5 Years: 0 Months : 0 Days",   "1826", "6.5",  "30 Apr 2029", "4000",    "₹1,300",   "65"),
    @("Quarterly Payout", "1500000", "5 Years: 0 Months : 0 Days",   "1826", "6.5",  "30 Apr 2029", "1500000", "₹487,500", "24,375"),
    @("Quarterly Payout", "2000",    "3 Years: 0 Months : 0 Days",   "1095", "7.25", "30 Apr 2027", "2000",    "₹435",     "36"),
    @("Quarterly Payout", "1000000", "6 Years: 0 Months : 0 Days",   "2191", "6.5",  "30 Apr 2030", "1000000", "₹390,000", "16,250")
)

$startRow = 32
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $values = $newRows[$i]
    for ($c = 1; $c -le $values.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        # Force text storage (matching the rest of the sheet, where every
        # value - including the numeric-looking ones - is a shared string,
        # not a number) without leaving a stray number-format style behind.
        $cell.Value = "'" + $values[$c - 1]
        $cell.Style = "Normal"
    }
}
